$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "project" hours for week 15 (row 16) from 34 to 36
$ws.Range("F16").Value = 36

# Move the active selection to F17, matching the saved cursor position
$ws.Activate()
$ws.Range("F17").Select()
